$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "binden",
    "reizen",
    "quälen",
    "sinken",
    "siegen",
    "seufzen",
    "meinen",
    "freuen",
    "filmen",
    "schütteln",
    "arten",
    "geben",
    "kratzen",
    "treffen",
    "wirken",
    "fühlen",
    "schrecken",
    "kümmern",
    "erben",
    "stoppen",
    "dienen",
    "stammen",
    "weichen",
    "fällen",
    "kaufen",
    "räumen",
    "saufen",
    "ärgern",
    "klingeln",
    "biegen",
    "regnen",
    "klettern",
    "krachen",
    "warnen",
    "zünden",
    "fließen",
    "dauern",
    "grüßen",
    "stürzen",
    "läuten",
    "machen",
    "ändern",
    "mauern",
    "werden",
    "sorgen",
    "merken",
    "süßen",
    "loben",
    "liefern",
    "decken",
    "flüchten",
    "ehren",
    "heben",
    "leisten",
    "sterben",
    "wählen",
    "pflanzen",
    "enden",
    "folgen",
    "schaden",
    "träumen",
    "ruhen",
    "fangen",
    "trennen",
    "bauen",
    "kehren",
    "trauen",
    "hören",
    "stellen",
    "sprengen",
    "fordern",
    "drehen",
    "backen",
    "stehlen",
    "bitten",
    "sperren",
    "rufen",
    "schicken",
    "scheitern",
    "zählen",
    "fesseln",
    "wehen",
    "achten",
    "schämen",
    "schlucken",
    "warten",
    "helfen",
    "schulden",
    "feiern",
    "kosten",
    "bellen",
    "laufen",
    "wüten",
    "liegen",
    "pissen",
    "wenden",
    "äußern",
    "herzen",
    "zeigen",
    "gelten",
    "rasen",
    "knurren",
    "schneiden",
    "töten",
    "bluten",
    "schenken",
    "zögern",
    "schwimmen",
    "bleiben",
    "boxen",
    "pfeifen",
    "lesen",
    "jubeln",
    "schwingen",
    "drucken",
    "wehtun",
    "füttern",
    "brauchen",
    "fischen",
    "bergen",
    "werfen",
    "messen",
    "fallen",
    "klappen",
    "mögen",
    "wetten",
    "starren",
    "malen",
    "spielen",
    "schwächen",
    "führen",
    "locken",
    "sichern",
    "danken",
    "bieten",
    "zielen",
    "wiegen",
    "dringen",
    "scheinen",
    "altern",
    "greifen",
    "graben",
    "scheiden",
    "stärken",
    "fassen",
    "flehen",
    "leiden",
    "runden",
    "wahren",
    "schließen",
    "suchen",
    "wachsen",
    "heulen",
    "klingen",
    "jagen",
    "münzen",
    "löschen",
    "rauben",
    "spinnen",
    "schmecken",
    "lockern",
    "füllen",
    "nerven",
    "tollen",
    "stecken",
    "gnaden",
    "opfern",
    "garen",
    "irren",
    "schreiten",
    "platzen",
    "flüstern",
    "zeugen",
    "planen",
    "ächzen",
    "kichern",
    "spüren",
    "lügen",
    "tropfen",
    "stoßen",
    "schwören",
    "heilen",
    "betteln",
    "streichen",
    "formen",
    "treiben",
    "wohnen",
    "gründen",
    "knarren",
    "hauen",
    "fahren",
    "wundern"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}
